# fix: atualiza planilha de gastos do t3
#
# The "T3 Microcontroladores" expense sheet is replaced by a (much
# shorter) T3 expense sheet: the list of components is replaced by a
# single "Amazon RDS (instância micro)" line item, and the total is
# recomputed as the sum of the (now single) cost cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old component rows (D9:E13 - Raspberry Pi 4GB Ram, Jumpers,
# Protoboard, Cartão SD, and the old pre-computed "Custo Total" row),
# shifting the rest of the sheet up. After this, row 9 is free again for
# the new "Custo Total" row.
$ws.Rows("9:13").Delete()

# Replace the remaining component row (D8:E8) with the new line item.
$ws.Range("D8").Value = "Amazon RDS (instância micro)"
$ws.Range("E8").Value = 45.53

# Recompute "Custo Total" (D9/E9) as a live sum over the cost column.
$ws.Range("E9").Formula = "=SUM(E8)"

# Match the author's final selection state.
[void]$ws.Range("E10").Select()
